$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Apply right-alignment (new style index 2) to the whole numeric data
#    block C2:J10 first, so that style becomes cellXfs index 2 before
#    any other style gets created.
# ---------------------------------------------------------------------
$ws.Range("C2:J10").HorizontalAlignment = -4152

# ---------------------------------------------------------------------
# 2) Fill in the Tesla rows (5=2023, 6=2022, 7=2021).
#    Numeric cells are plain values; a handful of operating-expense
#    cells were originally typed as text (trailing space keeps them as
#    text in the source workbook), so those are produced via a text
#    formula + paste-values trick so they land as shared-string text
#    cells (t="s") without dragging in a stray "@" number format.
# ---------------------------------------------------------------------

# Row 5 - Tesla 2023
$ws.Range("C5").Value = 96773
$ws.Range("D5").Value = 14974
$ws.Range("E5").Value = 106618
$ws.Range("F5").Value = 43009
$ws.Range("H5").Value = 0

# Row 6 - Tesla 2022
$ws.Range("D6").Value = 12587
$ws.Range("E6").Value = 82338
$ws.Range("F6").Value = 36440
$ws.Range("H6").Value = 0

# Row 7 - Tesla 2021
$ws.Range("D7").Value = 5644
$ws.Range("E7").Value = 62131
$ws.Range("F7").Value = 30548
$ws.Range("H7").Value = 0

# Text-typed numbers (column C for rows 6-7, columns G & I for rows 5-7)
# Built in the same order the original workbook's shared-string table
# uses: C6, C7, G5, G6, G7, I5, I6, I7.
$ws.Range("Z1").Formula = "=""81462 """
$ws.Range("Z1").Copy()
$ws.Range("C6").PasteSpecial(-4163)

$ws.Range("Z1").Formula = "=""53823 """
$ws.Range("Z1").Copy()
$ws.Range("C7").PasteSpecial(-4163)

$ws.Range("Z1").Formula = "=""3969 """
$ws.Range("Z1").Copy()
$ws.Range("G5").PasteSpecial(-4163)

$ws.Range("Z1").Formula = "=""3075 """
$ws.Range("Z1").Copy()
$ws.Range("G6").PasteSpecial(-4163)

$ws.Range("Z1").Formula = "=""2593 """
$ws.Range("Z1").Copy()
$ws.Range("G7").PasteSpecial(-4163)

$ws.Range("Z1").Formula = "=""4800 """
$ws.Range("Z1").Copy()
$ws.Range("I5").PasteSpecial(-4163)

$ws.Range("Z1").Formula = "=""3946 """
$ws.Range("Z1").Copy()
$ws.Range("I6").PasteSpecial(-4163)

$ws.Range("Z1").Formula = "=""4517 """
$ws.Range("Z1").Copy()
$ws.Range("I7").PasteSpecial(-4163)

$ws.Range("Z1").Clear()

# J5:J7 stay blank but still carry the new style (set below).

# Re-apply the alignment across the rows we just touched so every
# populated / blank cell in C5:J7 carries the shared style index 2.
$ws.Range("C5:J7").HorizontalAlignment = -4152

# ---------------------------------------------------------------------
# 3) Rows 8-10 (Apple) gain blank, styled cells C:J (no values).
# ---------------------------------------------------------------------
$ws.Range("C8:J10").HorizontalAlignment = -4152

# ---------------------------------------------------------------------
# 4) Selection moves from F16 to F10.
# ---------------------------------------------------------------------
$ws.Range("F10").Select()
